# The "Python Crash Course - Review" slide (originally slide #2) is moved
# later in the deck, to sit right before "recommended materials" (now slide #5,
# i.e. right after "Homework - part 2"). No slide content is edited - this is
# purely a slide-order change.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$s.MoveTo(5)
